# Populate the "Powerups" column (H) with values for the teams that were
# missing them, and refresh the corresponding "Score" column (J) values.
#
# Columns: A=Serial No, B=Team Name, C=Users, D=Phase Order,
#          E=Phase 1 Task Order, F=Phase 2 Task Order, G=Phase 3 Task Order,
#          H=Powerups, I=Credit Card No, J=Score

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Powerups = "1, 6, 7";    Score = 38  },
    @{ Row = 4;  Powerups = "6, 7, 1, 4"; Score = 279 },
    @{ Row = 5;  Powerups = "1, 6";       Score = 682 },
    @{ Row = 6;  Powerups = "1, 6, 2";    Score = 200 },
    @{ Row = 7;  Powerups = "5, 2, 6";    Score = 185 },
    @{ Row = 9;  Powerups = "1, 5";       Score = 539 },
    @{ Row = 10; Powerups = "7, 1, 5";    Score = 388 },
    @{ Row = 11; Powerups = "6, 7";       Score = 270 },
    @{ Row = 12; Powerups = "1, 6, 7";    Score = 644 },
    @{ Row = 13; Powerups = "2, 5";       Score = 156 },
    @{ Row = 19; Powerups = "7, 1";       Score = 337 },
    @{ Row = 20; Powerups = "9";          Score = 310 },
    @{ Row = 23; Powerups = "7, 5, 2";    Score = 23  },
    @{ Row = 27; Powerups = "7, 1";       Score = 521 }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, 8)
    if ($u.Powerups -match '^\d+$') {
        # Purely-numeric powerup lists (e.g. "9") must still be stored as
        # text, matching the source data's string column type. Force text
        # formatting before the write, then restore the default style so no
        # stray per-cell formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Powerups
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Powerups
    }
    $ws.Cells.Item($u.Row, 10).Value = $u.Score
}
